# The deck currently carries two theme parts:
#   ppt/theme/theme1.xml  -> "Integral" theme (Red Violet color scheme) - used by the slide master
#   ppt/theme/theme2.xml  -> "Office Theme" (Office color scheme)      - used by the notes master
#
# The target commit swaps the two themes' content, so the slide master ends up
# using the Office color scheme (theme1.xml becomes "Office Theme"-flavoured)
# and the notes master ends up with the Red Violet colors.
#
# The font scheme and format scheme (fills/lines/effects) are already
# byte-identical between the two themes, so the only real difference is the
# 12-color scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink). We apply
# that swap through the slide master's ColorScheme, which is the supported,
# round-trippable surface for editing theme colors via the PowerPoint object
# model.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

function Set-ThemeColor($scheme, $index, $r, $g, $b) {
    $color = $scheme.Colors($index)
    $color.RGB = $r + ($g * 256) + ($b * 65536)
}

# Office color scheme (was theme2.xml / notes master) -> now applied to the
# slide master's theme (theme1.xml), matching the new "Office Theme" colors.
Set-ThemeColor $colorScheme 1  0x00 0x00 0x00   # dk1
Set-ThemeColor $colorScheme 2  0xFF 0xFF 0xFF   # lt1
Set-ThemeColor $colorScheme 3  0x44 0x54 0x6A   # dk2
Set-ThemeColor $colorScheme 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeColor $colorScheme 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeColor $colorScheme 6  0xED 0x7D 0x31   # accent2
Set-ThemeColor $colorScheme 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeColor $colorScheme 8  0xFF 0xC0 0x00   # accent4
Set-ThemeColor $colorScheme 9  0x44 0x72 0xC4   # accent5
Set-ThemeColor $colorScheme 10 0x70 0xAD 0x47   # accent6
Set-ThemeColor $colorScheme 11 0x05 0x63 0xC1   # hlink
Set-ThemeColor $colorScheme 12 0x95 0x4F 0x72   # folHlink
